$d = $word.ActiveDocument

# Remember which paragraph is currently the last one in the document —
# the new content goes right after it.
$lastParaIndex = $d.Paragraphs.Count

# 1) Append the new "Under branch april-brn, add this line." paragraph at
#    the very end of the document. We build it from raw WordML via
#    InsertXML so we can place the w:proofErr spell-check markers (which
#    split the sentence into three runs) exactly like the target markup.
$endRange = $d.Content
$endRange.Collapse(0)  # wdCollapseEnd

$paragraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:r><w:t xml:space="preserve">Under branch </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>april-brn</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, add this line.</w:t></w:r>' +
    '</w:p>'
[void]$endRange.InsertXML($paragraphXml)

# 2) Insert a blank paragraph between the old last paragraph ("Updated for
#    a second commit.") and the paragraph we just appended.
$priorPara = $d.Paragraphs.Item($lastParaIndex)
$insertionPoint = $priorPara.Range
$insertionPoint.Collapse(0)  # wdCollapseEnd
$insertionPoint.InsertAfter("`r")
